$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 124800.04
$ws.Range("I15").Value = 124800.04
$ws.Range("K15").Value = 374400.12
$ws.Range("M15").Value = -374231.12
$ws.Range("H21").Value = 66250
$ws.Range("I21").Value = 75000
$ws.Range("J21").Value = 63333.332
$ws.Range("K21").Value = 75000
$ws.Range("L21").Value = 63333.332
$ws.Range("M21").Value = -74532
$ws.Range("N21").Value = -64269.332
$ws.Range("H23").Value = 66250
$ws.Range("I23").Value = 75000
$ws.Range("J23").Value = 63333.332
$ws.Range("K23").Value = 75000
$ws.Range("L23").Value = 63333.332
$ws.Range("M23").Value = -74766
$ws.Range("N23").Value = -63801.332
$ws.Range("H33").Value = 479.94736
$ws.Range("J33").Value = 199.5
$ws.Range("L33").Value = 199.5
$ws.Range("N33").Value = -657.5
$ws.Range("H112").Value = 3664.9688
$ws.Range("J112").Value = 3170.2903
$ws.Range("L112").Value = 9510.8709
$ws.Range("N112").Value = -11726.8709
$ws.Range("H125").Value = 2464.75
$ws.Range("I125").Value = 3410
$ws.Range("K125").Value = 30690
$ws.Range("M125").Value = -28230
$ws.Range("H137").Value = 2161.9268
$ws.Range("I137").Value = 770.3200000000001
$ws.Range("K137").Value = 2310.96
$ws.Range("M137").Value = 239.04

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4753.4287
$ws.Range("I2").Value = 2000.5
$ws.Range("K2").Value = 2000.5
$ws.Range("M2").Value = -1887.5
$ws.Range("H32").Value = 17371710
$ws.Range("I32").Value = 21746042
$ws.Range("K32").Value = 21746042
$ws.Range("M32").Value = -21745755
$ws.Range("H34").Value = 5024999.5
$ws.Range("I34").Value = 5024999.5
$ws.Range("K34").Value = 5024999.5
$ws.Range("M34").Value = -5024728.5
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H88").Value = 4156.1055
$ws.Range("J88").Value = 4026
$ws.Range("L88").Value = 4026
$ws.Range("N88").Value = -4838
$ws.Range("H91").Value = 4156.1055
$ws.Range("J91").Value = 4026
$ws.Range("L91").Value = 4026
$ws.Range("N91").Value = -6834
$ws.Range("H116").Value = 4753.4287
$ws.Range("I116").Value = 2000.5
$ws.Range("K116").Value = 2000.5
$ws.Range("M116").Value = 293.5
$ws.Range("H122").Value = 4904.1724
$ws.Range("I122").Value = 4356.737
$ws.Range("J122").Value = 5944.3
$ws.Range("K122").Value = 13070.211
$ws.Range("L122").Value = 17832.9
$ws.Range("M122").Value = -10620.211
$ws.Range("N122").Value = -22732.9

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4753.4287
$ws.Range("I3").Value = 2000.5
$ws.Range("K3").Value = 2000.5
$ws.Range("M3").Value = -1886.5
$ws.Range("H8").Value = 5004999.5
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 5004999.5
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 5004999.5
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -5005279.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 18280
$ws.Range("I93").Value = 18280
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 18280
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -16408
$ws.Range("N93").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 4983.1665
$ws.Range("I17").Value = 5724.75
$ws.Range("J17").Value = 3500
$ws.Range("K17").Value = 17174.25
$ws.Range("L17").Value = 10500
$ws.Range("M17").Value = -17005.25
$ws.Range("N17").Value = -10838
$ws.Range("H33").Value = 394
$ws.Range("I33").Value = 210.75
$ws.Range("J33").Value = 498.7143
$ws.Range("K33").Value = 1264.5
$ws.Range("L33").Value = 2992.2858
$ws.Range("M33").Value = -981.5
$ws.Range("N33").Value = -3558.2858
$ws.Range("H34").Value = 2372.7778
$ws.Range("I34").Value = 267.25
$ws.Range("J34").Value = 6583.8335
$ws.Range("K34").Value = 801.75
$ws.Range("L34").Value = 19751.5005
$ws.Range("M34").Value = -717.75
$ws.Range("N34").Value = -19919.5005
$ws.Range("H50").Value = 4050.2
$ws.Range("J50").Value = 4124
$ws.Range("L50").Value = 12372
$ws.Range("N50").Value = -13334
$ws.Range("H53").Value = 4050.2
$ws.Range("J53").Value = 4124
$ws.Range("L53").Value = 12372
$ws.Range("N53").Value = -13334
$ws.Range("H97").Value = 726.8333
$ws.Range("I97").Value = 254.11111
$ws.Range("K97").Value = 762.3333299999999
$ws.Range("M97").Value = -266.3333299999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 16953708
$ws.Range("I80").Value = 52633930
$ws.Range("K80").Value = 52633930
$ws.Range("M80").Value = -52632932
$ws.Range("H83").Value = 16953708
$ws.Range("I83").Value = 52633930
$ws.Range("K83").Value = 263169650
$ws.Range("M83").Value = -263164658
$ws.Range("H97").Value = 1223.52
$ws.Range("I97").Value = 1314.2778
$ws.Range("K97").Value = 1314.2778
$ws.Range("M97").Value = -818.2778000000001
$ws.Range("H102").Value = 6211.9487
$ws.Range("I102").Value = 5026.409
$ws.Range("J102").Value = 7746.1763
$ws.Range("K102").Value = 5026.409
$ws.Range("L102").Value = 7746.1763
$ws.Range("M102").Value = -3404.409
$ws.Range("N102").Value = -10990.1763

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 586.65515
$ws.Range("I22").Value = 531.2941
$ws.Range("J22").Value = 665.0833
$ws.Range("K22").Value = 531.2941
$ws.Range("L22").Value = 665.0833
$ws.Range("M22").Value = -236.2941
$ws.Range("N22").Value = -1255.0833
$ws.Range("H27").Value = 586.65515
$ws.Range("I27").Value = 531.2941
$ws.Range("J27").Value = 665.0833
$ws.Range("K27").Value = 531.2941
$ws.Range("L27").Value = 665.0833
$ws.Range("M27").Value = -424.2941
$ws.Range("N27").Value = -879.0833

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1509.88
$ws.Range("I100").Value = 1421
$ws.Range("J100").Value = 1791.3334
$ws.Range("K100").Value = 2842
$ws.Range("L100").Value = 3582.6668
$ws.Range("M100").Value = -2301
$ws.Range("N100").Value = -4664.6668
